$wb = $excel.ActiveWorkbook

$wsPreguntas = $wb.Worksheets.Item("PREGUNTAS")
$wsRespuestas = $wb.Worksheets.Item("RESPUESTAS")

# --- PREGUNTAS sheet: add a new "NIVEL" column (C) ------------------------

# Header cell: same look & feel as the existing "PREGUNTA"/"COMENTARIO" header cells
$wsPreguntas.Range("C1").Value = "NIVEL"
$wsPreguntas.Range("A1").Copy() | Out-Null
$wsPreguntas.Range("C1").PasteSpecial(-4122) | Out-Null

# Data cells: numeric "level" values, centered like the rest of the sheet
$wsPreguntas.Range("C2").Value = 1
$wsPreguntas.Range("C3").Value = 2
$wsPreguntas.Range("C2:C3").HorizontalAlignment = -4108

# --- Selection / active sheet ---------------------------------------------
# The edited file now opens on PREGUNTAS (first tab) instead of RESPUESTAS.
$wsPreguntas.Activate() | Out-Null
$wsPreguntas.Range("C4").Select() | Out-Null
